$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table with a new "2023" column (K), copying the formatting
# of the preceding "2022" column (J) so the new column looks identical.
$ws.Range("J3:J6").Copy()
$ws.Range("K3:K6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header for the new year column.
$ws.Range("K3").Value = 2023

# Carry the 2022 figures forward into 2023 (placeholder values, same as
# column J) for each of the three data rows.
$ws.Range("K4").Value = $ws.Range("J4").Value2
$ws.Range("K5").Value = $ws.Range("J5").Value2
$ws.Range("K6").Value = $ws.Range("J6").Value2
